$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 9 & 10: "Garra nana" / "Garra rufa" become "Garra jordanica" / "Garra nana" ---
# The numeric data (C:R) of row 9 and row 10 swap places; row 9 gets a brand new
# species name "Garra jordanica" while row 10 inherits row 9's old name "Garra nana".

# Remember the old name that used to be in B9 before we overwrite anything.
$oldB9 = $ws.Range("B9").Value()

# Swap the numeric blocks C9:R9 <-> C10:R10 using a temporary holding row (row 40).
$ws.Range("C9:R9").Copy($ws.Range("C40:R40"))
$ws.Range("C10:R10").Copy($ws.Range("C9:R9"))
$ws.Range("C40:R40").Copy($ws.Range("C10:R10"))
$ws.Range("C40:R40").Clear()

# Update the species names.
$ws.Range("B9").Value = "Garra jordanica"
$ws.Range("B10").Value = $oldB9

# --- Row 14: "Mugil liza" -> "Mugil" (data untouched) ---
$ws.Range("B14").Value = "Mugil"

# --- Rows 17-20: cyclic rotation ---
# Data (C:R) of row 20 moves to row 17 (and row 17 gets a brand-new name "Coptodon zillii"),
# while rows 17,18,19 (name + data) shift down into 18,19,20.

# Names shift down first (remember originals before overwriting).
$nameB17 = $ws.Range("B17").Value()
$nameB18 = $ws.Range("B18").Value()
$nameB19 = $ws.Range("B19").Value()

$ws.Range("B20").Value = $nameB19
$ws.Range("B19").Value = $nameB18
$ws.Range("B18").Value = $nameB17
$ws.Range("B17").Value = "Coptodon zillii"

# Rotate the numeric data blocks C17:R17 -> C18:R18 -> C19:R19 -> C20:R20 -> C17:R17
$ws.Range("C20:R20").Copy($ws.Range("C40:R40"))
$ws.Range("C19:R19").Copy($ws.Range("C20:R20"))
$ws.Range("C18:R18").Copy($ws.Range("C19:R19"))
$ws.Range("C17:R17").Copy($ws.Range("C18:R18"))
$ws.Range("C40:R40").Copy($ws.Range("C17:R17"))
$ws.Range("C40:R40").Clear()
